$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (rnn_006): comment was missing the Tsim=1200 note - fix it
$ws.Range("I15").Value = "random uniform noise added to hidden states, excluding the first timestep, dataset normalized (featurewise) with sample mean and std, Tsim=1200"
$ws.Range("I15").Orientation = 0

# Row 21 (new rnn_009 entry) - copy formatting from row 19 first, then fill values
$ws.Range("A19:I19").Copy($ws.Range("A21:I21"))
$ws.Range("A21").Value = "colab"
$ws.Range("B21").Value = "rnn_009"
$ws.Range("C21").Value = "[64]"
$ws.Range("D21").Value = "GRU"
$ws.Range("E21").Value = "LR sigmoid (warmup 20, expected 50)"
$ws.Range("F21").Value = 90.0
$ws.Range("G21").Value = 90.0
$ws.Range("H21").Value = 0.1
$ws.Range("I21").Value = "random uniform noise (stddev=1e-4) added to hidden states, excluding the first timestep, dataset normalized (featurewise) with sample mean and std, Tsim=1000"
$ws.Range("B21").Orientation = 0
$ws.Range("I21").Orientation = 0
